# Generate Report for handback
# Marks the c4163379-... handback row as complete ("Handed back: in sync
# with en-US") on every sheet that shows it, fills in the Latest Target
# File / Latest Handback File hyperlinks and the Latest Handback DateTime
# for both locales.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$mdUrl      = "https://github.com/OpenLocalizationTest/oltest/blob/0615ce0a3282af785c99179610147afdd4a5a83b/e2e/c4163379-5b97-49f3-a079-0ab1916eeb61.md"
$zhHandbackUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8bfb06509ce88e67f36ca74b9df523d6c9aa8b8c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/c4163379-5b97-49f3-a079-0ab1916eeb61.9db323ca607701c0f96fd9f657f4c7c1417f9959.zh-cn.xlf"
$deHandbackUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/337c0165172fbbbf69f8fe07e5c7e5489ea8b71f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/c4163379-5b97-49f3-a079-0ab1916eeb61.9db323ca607701c0f96fd9f657f4c7c1417f9959.de-de.xlf"
$configUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/0615ce0a3282af785c99179610147afdd4a5a83b/.localization-config"

$mdName = "c4163379-5b97-49f3-a079-0ab1916eeb61.md"
$zhXlfName = "c4163379-5b97-49f3-a079-0ab1916eeb61.9db323ca607701c0f96fd9f657f4c7c1417f9959.zh-cn.xlf"
$deXlfName = "c4163379-5b97-49f3-a079-0ab1916eeb61.9db323ca607701c0f96fd9f657f4c7c1417f9959.de-de.xlf"
$configName = ".localization-config"

# ---- Overview sheet: refresh the status shown for both locales ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("B2").Value = $newStatus
$wsZh.Range("G2").Value = "2016-01-14 06:58:15"

# Rebuild the hyperlinks in cell order (A2, C2, E2, F2, A3) so the two new
# links land in the middle of the collection, matching a freshly generated
# report.
$wsZh.Hyperlinks.Delete()
$wsZh.Range("E2").Value = $mdName
$wsZh.Range("F2").Value = $zhXlfName
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $mdUrl, "", "", $mdName)
$wsZh.Hyperlinks.Add($wsZh.Range("C2"), $zhHandbackUrl, "", "", $zhXlfName)
$wsZh.Hyperlinks.Add($wsZh.Range("E2"), $mdUrl, "", "", $mdName)
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), $zhHandbackUrl, "", "", $zhXlfName)
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $configUrl, "", "", $configName)

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("B2").Value = $newStatus
$wsDe.Range("G2").Value = "2016-01-14 06:58:54"

$wsDe.Hyperlinks.Delete()
$wsDe.Range("E2").Value = $mdName
$wsDe.Range("F2").Value = $deXlfName
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $mdUrl, "", "", $mdName)
$wsDe.Hyperlinks.Add($wsDe.Range("C2"), $deHandbackUrl, "", "", $deXlfName)
$wsDe.Hyperlinks.Add($wsDe.Range("E2"), $mdUrl, "", "", $mdName)
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), $deHandbackUrl, "", "", $deXlfName)
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $configUrl, "", "", $configName)
